$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data for 29/11
$ws.Range("A5").Value = "29/11"
$ws.Range("B5").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B5").Value = 1.0 / 24.0
$ws.Range("C5").Value = "Restructuring character + basic sword implementation + camera movement"
